$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.121.27'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.831.93'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4624'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3704'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07357'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8731'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07986'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.87'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.05%  '
$ws.Range('D13').Value = '1.794.95'
$ws.Range('E13').Value = '  -5.57%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.347'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.574'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '92.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008884'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.70'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '27.169.97'
$ws.Range('E21').Value = '  -1.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.145'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '2.061.60'
$ws.Range('E24').Value = '  -1.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.836'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.59'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.090'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.092'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.50'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08870'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.973'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7338'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.450'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.139'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.460'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.075'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01947'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05241'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.940'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.179'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5196'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.26%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1635'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.54%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8591'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -14.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.242'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4844'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.28'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.008'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.632'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06236'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.72%  '
